$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Settings Page" task to reflect the new feature work
$ws.Range("A2").Value = "Settings Page add company add user"

# Update "Time Taken" (column C) hours for the in-progress tasks
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 1

# Update the Status column fill colors (D2:D4) to reflect progress:
#   green  = 00B050 -> Done
#   orange = FFC000 -> Partially Done
#   red    = C00000 -> Not Started
$ws.Range("D2").Interior.Color = 49407   # Partially Done (orange)
$ws.Range("D3").Interior.Color = 5287936 # Done (green)
$ws.Range("D4").Interior.Color = 49407   # Partially Done (orange)

# Move the active selection
$ws.Range("C14").Select()
